$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updated values (Modified Reg iProctor TC's)
$ws.Range("A2").Value = "AqtnV892"
$ws.Range("B2").Value = 23073123
$ws.Range("C2").Value = "ifikrge80"
$ws.Range("D2").Value = "Z&2At7d$"
$ws.Range("F2").Value = "NKBJXgBO"
$ws.Range("G2").Value = "upIJ"
